# Update exported_companys.xlsx export:
#  - Row 2 (the "WayBack" brand row) loses its extra 6th column ("wayback")
#    and its "United States of America NSA" value moves from F2 into E2.
#  - Row 3's last cell ("demo-company") is cleared to an empty string.
#  - A new row 4 is appended that duplicates row 2's (new) 5-column content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture row 2's current values up front (read via the Value() accessor -
# plain property access on this COM shim doesn't invoke the getter).
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$d2 = $ws.Range("D2").Value()
$nsa = $ws.Range("F2").Value()

# Move "United States of America NSA" from F2 into E2, then drop F2.
$ws.Range("E2").Value = $nsa
$ws.Range("F2").Value = $null

# Row 3, column E ("demo-company") becomes an empty string.
$ws.Range("E3").Value = ""

# New row 4: duplicate of row 2's (new) 5-column content.
$ws.Range("A4").Value = $a2
$ws.Range("B4").Value = $b2
$ws.Range("C4").Value = $c2
$ws.Range("D4").Value = $d2
$ws.Range("E4").Value = $nsa
